$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.611.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.60%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  +7.20%  "

# Row 9
$ws.Range("E9").Value = "  +0.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.051.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.797.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.634.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.631"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("E17").Value = "  +2.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "254.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0791"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.23%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "

# Row 27
$ws.Range("E27").Value = "  +1.07%  "

# Row 28
$ws.Range("E28").Value = "  -0.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("E30").Value = "  +0.22%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.71%  "

# Row 32
$ws.Range("E32").Value = "  -0.50%  "

# Row 33
$ws.Range("E33").Value = "  -0.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.55%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.438.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.27%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "

# Row 38
$ws.Range("E38").Value = "  +2.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.86%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.915"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "

# Row 44
$ws.Range("E44").Value = "  +5.19%  "

# Row 45
$ws.Range("E45").Value = "  -1.56%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0494"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.950.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.79%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "

# Row 50
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.12%  "
